$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: copy formatting (date/time styles) from row 15, then set the new values
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("A16").Value = 41968
$ws.Range("B16").Value = 0.54166666666666663

# Row 17: copy formatting (date/time styles) from row 15, then set the new values
$ws.Range("A15:B15").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)
$ws.Range("A17").Value = 41976
$ws.Range("B17").Value = 0.5

$excel.CutCopyMode = 0

$ws.Range("B17").Select()
